$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '27.789.10'
$c.Style = "Normal"
$ws.Range('E2').Value = '  +1.35%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '1.650.31'
$c.Style = "Normal"
$ws.Range('E3').Value = '  -0.29%  '
$ws.Range('E4').Value = '  +0.11%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '213.56'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('E7').Value = '  +0.11%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '23.19'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -1.02%  '
$ws.Range('E9').Value = '  -0.21%  '
$ws.Range('E10').Value = '  +0.52%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0892'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -1.62%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '1.884.27'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -0.25%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '1.655.20'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('E15').Value = '  -0.21%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '64.53'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -1.47%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '27.772.09'
$c.Style = "Normal"
$ws.Range('E17').Value = '  +1.31%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '234.11'
$c.Style = "Normal"
$ws.Range('E18').Value = '  +2.51%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '7.71'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +3.84%  '
$ws.Range('E20').Value = '  -0.10%  '
$ws.Range('E21').Value = '  +0.12%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '4.32'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -0.54%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '10.15'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +8.13%  '
$ws.Range('E24').Value = '  -3.85%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '150.76'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +2.58%  '
$ws.Range('E26').Value = '  -1.06%  '
$ws.Range('E27').Value = '  -1.85%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '15.70'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('E30').Value = '  +0.38%  '
$ws.Range('E31').Value = '  -1.02%  '
$ws.Range('E32').Value = '  +0.68%  '
$ws.Range('E33').Value = '  +1.68%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '1.443.84'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +1.63%  '
$ws.Range('E35').Value = '  +1.93%  '
$ws.Range('E36').Value = '  -1.04%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.573'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +0.95%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.887'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -2.05%  '
$ws.Range('E39').Value = '  -0.22%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.880'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +11.62%  '
$ws.Range('E41').Value = '  -1.21%  '
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('E43').Value = '  +1.42%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '66.70'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +2.44%  '
$ws.Range('E45').Value = '  -0.84%  '
$ws.Range('E46').Value = '  +2.19%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '1.793.00'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('E48').Value = '  +4.39%  '
$ws.Range('E49').Value = '  -1.61%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.0₆0107'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +2.21%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.0998'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -1.16%  '
